# Edit word 'process' to 'proceed'
# The column_ref sheet holds the canonical header labels that every e2e_*
# sheet's row-1 pulls in via an array formula (=column_ref!A1:O1), so fixing
# the single source cell on column_ref propagates the corrected text to
# every dependent sheet automatically.

$wb = $excel.ActiveWorkbook

$columnRef = $wb.Worksheets.Item("column_ref")
$columnRef.Range("J1").Value = "proceedToCheckoutFrom"

# --- Recreate the recorded view state (selections / zoom / active tab) ---

$ws = $wb.Worksheets.Item("column_ref")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 70
[void]$ws.Range("M7").Select()

$ws = $wb.Worksheets.Item("e2e_a02")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 60
[void]$ws.Range("Q16").Select()

$ws = $wb.Worksheets.Item("e2e_a03")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 60
[void]$ws.Range("I2").Select()

$ws = $wb.Worksheets.Item("e2e_a04")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 60
[void]$ws.Range("P3").Select()

$ws = $wb.Worksheets.Item("e2e_a06")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 60
[void]$ws.Range("O21").Select()

$ws = $wb.Worksheets.Item("e2e_a07")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 60
[void]$ws.Range("O17").Select()

$ws = $wb.Worksheets.Item("e2e_e02")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 60
[void]$ws.Range("J4").Select()

# e2e_a01 is the sheet left active/selected when the workbook was saved.
$ws = $wb.Worksheets.Item("e2e_a01")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 60
[void]$ws.Range("B5").Select()
